$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(154, 2).Value = 57756
$ws.Cells.Item(154, 6).Value = -100
$ws.Cells.Item(154, 7).Value = -6644

$ws.Cells.Item(156, 2).Value = 53925
$ws.Cells.Item(156, 6).Value = 1
$ws.Cells.Item(156, 7).Value = 66.44

$ws.Cells.Item(176, 2).Value = 57552
$ws.Cells.Item(176, 5).Value = 136.86
$ws.Cells.Item(176, 6).Value = -5
$ws.Cells.Item(176, 7).Value = -603.45

$ws.Cells.Item(177, 2).Value = 64329
$ws.Cells.Item(177, 5).Value = 128.32
$ws.Cells.Item(177, 6).Value = 6
$ws.Cells.Item(177, 7).Value = 724.14

$ws.Cells.Item(271, 2).Value = 64973
$ws.Cells.Item(271, 5).Value = 35.4
$ws.Cells.Item(271, 6).Value = 150
$ws.Cells.Item(271, 7).Value = 4995

$ws.Cells.Item(272, 2).Value = 48706
$ws.Cells.Item(272, 5).Value = 39.8
$ws.Cells.Item(272, 6).Value = -144
$ws.Cells.Item(272, 7).Value = -4795.2

$ws.Cells.Item(308, 2).Value = 57077
$ws.Cells.Item(308, 4).Value = 93.08
$ws.Cells.Item(308, 5).Value = 111.2
$ws.Cells.Item(308, 6).Value = 1
$ws.Cells.Item(308, 7).Value = 93.08

$ws.Cells.Item(309, 2).Value = 61610
$ws.Cells.Item(309, 5).Value = 122.71
$ws.Cells.Item(309, 6).Value = -58
$ws.Cells.Item(309, 7).Value = -5957.18

$ws.Cells.Item(310, 2).Value = 63565
$ws.Cells.Item(310, 4).Value = 102.71
$ws.Cells.Item(310, 5).Value = 109.19
$ws.Cells.Item(310, 6).Value = 60
$ws.Cells.Item(310, 7).Value = 6162.6

$ws.Cells.Item(338, 2).Value = 55373
$ws.Cells.Item(338, 5).Value = 163.62
$ws.Cells.Item(338, 6).Value = -94
$ws.Cells.Item(338, 7).Value = -13562.32

$ws.Cells.Item(339, 2).Value = 63520
$ws.Cells.Item(339, 5).Value = 153.4
$ws.Cells.Item(339, 6).Value = 97
$ws.Cells.Item(339, 7).Value = 13995.16

$ws.Cells.Item(342, 2).Value = 57802
$ws.Cells.Item(342, 5).Value = 162.71
$ws.Cells.Item(342, 6).Value = -79
$ws.Cells.Item(342, 7).Value = -11334.92

$ws.Cells.Item(343, 2).Value = 63531
$ws.Cells.Item(343, 5).Value = 152.53
$ws.Cells.Item(343, 6).Value = 80
$ws.Cells.Item(343, 7).Value = 11478.4

$ws.Cells.Item(344, 2).Value = 63571
$ws.Cells.Item(344, 6).Value = 29
$ws.Cells.Item(344, 7).Value = 4160.92

$ws.Cells.Item(347, 2).Value = 63510
$ws.Cells.Item(347, 5).Value = 50.66
$ws.Cells.Item(347, 6).Value = 167
$ws.Cells.Item(347, 7).Value = 7955.88

$ws.Cells.Item(348, 2).Value = 55356
$ws.Cells.Item(348, 5).Value = 54.04
$ws.Cells.Item(348, 6).Value = -158
$ws.Cells.Item(348, 7).Value = -7527.12

$ws.Cells.Item(371, 2).Value = 63564
$ws.Cells.Item(371, 5).Value = 137.16
$ws.Cells.Item(371, 6).Value = 57
$ws.Cells.Item(371, 7).Value = 7353.57

$ws.Cells.Item(372, 2).Value = 61608
$ws.Cells.Item(372, 5).Value = 154.12
$ws.Cells.Item(372, 6).Value = -56
$ws.Cells.Item(372, 7).Value = -7224.56

$ws.Cells.Item(374, 2).Value = 60325
$ws.Cells.Item(374, 5).Value = 151.57
$ws.Cells.Item(374, 6).Value = -102
$ws.Cells.Item(374, 7).Value = -12939.72

$ws.Cells.Item(375, 2).Value = 63560
$ws.Cells.Item(375, 5).Value = 134.87
$ws.Cells.Item(375, 6).Value = 104
$ws.Cells.Item(375, 7).Value = 13193.44

$ws.Cells.Item(381, 2).Value = 57817
$ws.Cells.Item(381, 6).Value = 3
$ws.Cells.Item(381, 7).Value = 239.43

$ws.Cells.Item(382, 2).Value = 62865
$ws.Cells.Item(382, 6).Value = 151
$ws.Cells.Item(382, 7).Value = 12051.31

$ws.Cells.Item(392, 2).Value = 62933
$ws.Cells.Item(392, 6).Value = 146
$ws.Cells.Item(392, 7).Value = 8632.98

$ws.Cells.Item(393, 2).Value = 57835
$ws.Cells.Item(393, 6).Value = 1
$ws.Cells.Item(393, 7).Value = 59.13

$ws.Cells.Item(423, 2).Value = 53082
$ws.Cells.Item(423, 3).Value = "HUL-VIM BAR MULTIPACK FW 4X200G"
$ws.Cells.Item(423, 6).Value = 1
$ws.Cells.Item(423, 7).Value = 59.47

$ws.Cells.Item(424, 2).Value = 63102
$ws.Cells.Item(424, 3).Value = "HUL-Vim Bar Multipack Fw 4X200G"
$ws.Cells.Item(424, 6).Value = 36
$ws.Cells.Item(424, 7).Value = 2140.92

$ws.Cells.Item(575, 2).Value = 65066
$ws.Cells.Item(575, 5).Value = 13.61
$ws.Cells.Item(575, 6).Value = 313
$ws.Cells.Item(575, 7).Value = 4009.53

$ws.Cells.Item(576, 2).Value = 53263
$ws.Cells.Item(576, 5).Value = 15.29
$ws.Cells.Item(576, 6).Value = -309
$ws.Cells.Item(576, 7).Value = -3958.29

$ws.Cells.Item(578, 2).Value = 64915
$ws.Cells.Item(578, 5).Value = 20.98
$ws.Cells.Item(578, 6).Value = 40
$ws.Cells.Item(578, 7).Value = 789.2

$ws.Cells.Item(579, 2).Value = 45695
$ws.Cells.Item(579, 5).Value = 23.58
$ws.Cells.Item(579, 6).Value = -36
$ws.Cells.Item(579, 7).Value = -710.28

$ws.Cells.Item(582, 2).Value = 64922
$ws.Cells.Item(582, 5).Value = 20.98
$ws.Cells.Item(582, 6).Value = 207
$ws.Cells.Item(582, 7).Value = 4084.11

$ws.Cells.Item(583, 2).Value = 45706
$ws.Cells.Item(583, 5).Value = 23.58
$ws.Cells.Item(583, 6).Value = -202
$ws.Cells.Item(583, 7).Value = -3985.46

$ws.Cells.Item(585, 2).Value = 45718
$ws.Cells.Item(585, 5).Value = 19.38
$ws.Cells.Item(585, 6).Value = -294
$ws.Cells.Item(585, 7).Value = -4768.68

$ws.Cells.Item(586, 2).Value = 64927
$ws.Cells.Item(586, 5).Value = 17.26
$ws.Cells.Item(586, 6).Value = 295
$ws.Cells.Item(586, 7).Value = 4784.9

$ws.Cells.Item(591, 2).Value = 64925
$ws.Cells.Item(591, 5).Value = 13.97
$ws.Cells.Item(591, 6).Value = 302
$ws.Cells.Item(591, 7).Value = 3971.3

$ws.Cells.Item(592, 2).Value = 45709
$ws.Cells.Item(592, 5).Value = 15.69
$ws.Cells.Item(592, 6).Value = -300
$ws.Cells.Item(592, 7).Value = -3945

$ws.Cells.Item(596, 2).Value = 53595
$ws.Cells.Item(596, 5).Value = 17.61
$ws.Cells.Item(596, 6).Value = -335
$ws.Cells.Item(596, 7).Value = -4934.55

$ws.Cells.Item(597, 2).Value = 65067
$ws.Cells.Item(597, 5).Value = 15.65
$ws.Cells.Item(597, 6).Value = 338
$ws.Cells.Item(597, 7).Value = 4978.74

$ws.Cells.Item(679, 2).Value = 53319
$ws.Cells.Item(679, 5).Value = 310.64
$ws.Cells.Item(679, 6).Value = -6
$ws.Cells.Item(679, 7).Value = -1643.52

$ws.Cells.Item(680, 2).Value = 64810
$ws.Cells.Item(680, 5).Value = 291.22
$ws.Cells.Item(680, 6).Value = 7
$ws.Cells.Item(680, 7).Value = 1917.44

$ws.Cells.Item(701, 2).Value = 64833
$ws.Cells.Item(701, 5).Value = 34.9
$ws.Cells.Item(701, 6).Value = 99
$ws.Cells.Item(701, 7).Value = 3250.17

$ws.Cells.Item(702, 2).Value = 60025
$ws.Cells.Item(702, 5).Value = 37.22
$ws.Cells.Item(702, 6).Value = -98
$ws.Cells.Item(702, 7).Value = -3217.34

$ws.Cells.Item(712, 2).Value = 64830
$ws.Cells.Item(712, 5).Value = 34.9
$ws.Cells.Item(712, 6).Value = 117
$ws.Cells.Item(712, 7).Value = 3841.11

$ws.Cells.Item(713, 2).Value = 60022
$ws.Cells.Item(713, 5).Value = 37.22
$ws.Cells.Item(713, 6).Value = -113
$ws.Cells.Item(713, 7).Value = -3709.79
